# Adds a new "2022-Q1" sheet of per-fund holdings (reusing the current
# "总计" sheet's slot) and rebuilds the "总计" (grand-total) summary sheet
# with a new leading row for 2022-Q1, pushing the rest down.

function Set-TextCell($cell, $text) {
    # Forces Excel to store the value as text instead of silently parsing
    # numeric-looking strings ("14.75", "006679", ...) into numbers, then
    # drops the temporary text number-format again so the cell is left
    # with no explicit style (matching the surrounding unstyled data
    # cells).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# The existing last sheet is "总计" (sheetId 4). We duplicate it first so
# the duplicate (placed right after it) can become the new "总计", while
# the original slot gets renamed + refilled to become "2022-Q1" - this
# mirrors how the sheetIds/rIds shift in the target workbook.
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalSheet.Copy($null, $totalSheet)
$newTotalSheet = $wb.Worksheets.Item($totalSheet.Index + 1)

# ---------------------------------------------------------------------
# 1) Extend formatting on the original sheet (soon to be "2022-Q1") so it
#    has 8 columns (A-H) and 5 rows, matching the header/index styling
#    already used for B1:D1 and A2:A4.
# ---------------------------------------------------------------------
$totalSheet.Range("B1:D1").Copy()
$totalSheet.Range("E1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Fill in the "2022-Q1" per-fund holding data.
# ---------------------------------------------------------------------
$totalSheet.Range("B1").Value = "基金代码"
$totalSheet.Range("C1").Value = "基金名称"
$totalSheet.Range("D1").Value = "基金规模"
$totalSheet.Range("E1").Value = "股票总仓位"
$totalSheet.Range("F1").Value = "仓位占比"
$totalSheet.Range("G1").Value = "持有市值(亿元)"
$totalSheet.Range("H1").Value = "仓位排名"

$totalSheet.Range("A2").Value = 0
Set-TextCell $totalSheet.Range("B2") "006679"
$totalSheet.Range("C2").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇A"
Set-TextCell $totalSheet.Range("D2") "14.75"
Set-TextCell $totalSheet.Range("E2") "83.19"
Set-TextCell $totalSheet.Range("F2") "6.34"
Set-TextCell $totalSheet.Range("G2") "0.9352"
$totalSheet.Range("H2").Value = 3

$totalSheet.Range("A3").Value = 1
Set-TextCell $totalSheet.Range("B3") "162719"
$totalSheet.Range("C3").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）A"
Set-TextCell $totalSheet.Range("D3") "14.75"
Set-TextCell $totalSheet.Range("E3") "83.19"
Set-TextCell $totalSheet.Range("F3") "6.34"
Set-TextCell $totalSheet.Range("G3") "0.9352"
$totalSheet.Range("H3").Value = 3

$totalSheet.Range("A4").Value = 2
Set-TextCell $totalSheet.Range("B4") "006680"
$totalSheet.Range("C4").Value = "广发道琼斯美国石油开发与生产指数证券投资基金(QDII-LOF) 美元现汇C"
Set-TextCell $totalSheet.Range("D4") "4.73"
Set-TextCell $totalSheet.Range("E4") "83.19"
Set-TextCell $totalSheet.Range("F4") "6.34"
Set-TextCell $totalSheet.Range("G4") "0.2999"
$totalSheet.Range("H4").Value = 3

$totalSheet.Range("A5").Value = 3
Set-TextCell $totalSheet.Range("B5") "004243"
$totalSheet.Range("C5").Value = "广发道琼斯美国石油开发与生产指数（QDII-LOF）C"
Set-TextCell $totalSheet.Range("D5") "4.73"
Set-TextCell $totalSheet.Range("E5") "83.19"
Set-TextCell $totalSheet.Range("F5") "6.34"
Set-TextCell $totalSheet.Range("G5") "0.2999"
$totalSheet.Range("H5").Value = 3

$totalSheet.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 3) Rebuild the (duplicated) "总计" sheet: same 3-column layout as
#    before, with a new 2022-Q1 row inserted at the top and the rest of
#    the quarters shifted down a row. The duplicate only has 4 data
#    rows (A1:D4), so extend the index-column styling down to row 5
#    before writing to it.
# ---------------------------------------------------------------------
$newTotalSheet.Range("A2").Copy()
$newTotalSheet.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newTotalSheet.Range("B1").Value = "日期"
$newTotalSheet.Range("C1").Value = "持有数量(只)"
$newTotalSheet.Range("D1").Value = "持有市值(亿元)"

$newTotalSheet.Range("A2").Value = 0
$newTotalSheet.Range("B2").Value = "2022-Q1"
$newTotalSheet.Range("C2").Value = 4
$newTotalSheet.Range("D2").Value = 2.47

$newTotalSheet.Range("A3").Value = 1
$newTotalSheet.Range("B3").Value = "2021-Q4"
$newTotalSheet.Range("C3").Value = 4
$newTotalSheet.Range("D3").Value = 0.75

$newTotalSheet.Range("A4").Value = 2
$newTotalSheet.Range("B4").Value = "2021-Q3"
$newTotalSheet.Range("C4").Value = 4
$newTotalSheet.Range("D4").Value = 0.79

$newTotalSheet.Range("A5").Value = 3
$newTotalSheet.Range("B5").Value = "2021-Q2"
$newTotalSheet.Range("C5").Value = 4
$newTotalSheet.Range("D5").Value = 0.5600000000000001

$newTotalSheet.Name = "总计"
